# Applies the ALASKA_2016.xlsx data-cleaning fixes:
#  - rename header columns to snake_case field names
#  - normalize "de" -> "De" in a handful of place names
#  - remove the trailing metadata/footer rows (52-56 and 476-480),
#    which also shrinks the used range back down to A1:D50

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two footer blocks first, starting from the bottom so that
# deleting the lower block doesn't shift the row numbers of the upper one.
$ws.Rows("476:480").Delete() | Out-Null
$ws.Rows("52:56").Delete() | Out-Null

# Rename header row (row 1) to the new snake_case field names.
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Capitalize "De" in a few place names.
$ws.Range("A12").Value = "Ciudad De México"
$ws.Range("A16").Value = "Estado De México"
$ws.Range("B21").Value = "Pachuca De Soto"
$ws.Range("B48").Value = "Tlaltenango De Sánchez Román"
